$wb = $excel.ActiveWorkbook

# --- Sheet1: rename Migraatio -> Migraatiot ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Migraatiot"

# --- Add new sheet "Tilat" right after the first sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tilat"

# --- Populate Tilat sheet (status list) ---
$ws2.Range("A1").Value = "Tila"
$ws2.Range("A2").Value = "SUUNNITTELU"
$ws2.Range("A3").Value = "NAHTAVILLAOLO"
$ws2.Range("A4").Value = "HYVAKSYMISPAATOS"
$ws2.Range("A5").Value = "JATKOPAATOS1"
$ws2.Range("A6").Value = "JATKOPAATOS2"

# --- Update Migraatiot row 2 sample data (order matters for shared-string layout) ---
$ws1.Range("B2").Value = "JATKOPAATOS2"
$ws1.Range("C2").Value = "TESTI-ASIA-123"
$ws1.Range("A2").Value = "1.2.246.578.5.1.2293640800.1682339657"
$ws1.Range("E2").Value = "TESTI-ASIA-243"

# Match C2/E2 formatting with the rest of the data row (B2)
$ws1.Range("B2").Copy() | Out-Null
$ws1.Range("C2").PasteSpecial(-4122) | Out-Null
$ws1.Range("E2").PasteSpecial(-4122) | Out-Null

# Dates for the approval / continuation decision columns
$ws1.Range("D2").NumberFormat = "yyyy-mm-dd;@"
$ws1.Range("D2").Value = 44611
$ws1.Range("F2").NumberFormat = "yyyy-mm-dd;@"
$ws1.Range("F2").Value = 44606

# --- Clear the now-unused sample cells in rows 3:10 and reset D/F date formatting ---
$ws1.Range("A3:F10").Clear()
$ws1.Range("D3:D13").NumberFormat = "yyyy-mm-dd;@"
$ws1.Range("F3:F13").NumberFormat = "yyyy-mm-dd;@"

# --- Data validation: Tila column now pulls from the Tilat table ---
$ws1.Range("B1").Validation.Delete()
$ws1.Range("B2:B13").Validation.Add(3, 1, 1, 'INDIRECT("Tilat[Tila]")')

# --- Turn the two ranges into real Excel Tables ---
$lo1 = $ws1.ListObjects.Add(1, $ws1.Range("A1:F13"), $null, 1)
$lo1.Name = "Migraatiot"
$lo1.TableStyle = "TableStyleLight15"

$lo2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:A6"), $null, 1)
$lo2.Name = "Tilat"

# --- Selection parity with the authored workbook ---
$ws1.Range("B10").Select() | Out-Null
